# meeting 01.07. todo update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("to dos")

# Row 7: new task cell in column C (was empty before)
$ws.Range("C7").Value = "UND suche eine Bildungswissenschaftliche Studie in der Tobii Glasses benutzt wurden, lege"

# Row 13: C13 replaced ("ziehe aus den drei Journals..." -> "Tobii Glasses Paper irgendwoher")
$ws.Range("C13").Value = "Tobii Glasses Paper irgendwoher"

# Row 13: F13 replaced ("GLMM-Artikel lesen + verstehen" -> "LMM-Artikel lesen + verstehen")
$ws.Range("F13").Value = "LMM-Artikel lesen + verstehen"

# Row 14: D14 replaced ("x" -> "Methodsection: Brille beschreiben")
$ws.Range("D14").Value = "Methodsection: Brille beschreiben"

# Row 15: D15 replaced ("x" -> "Eineleitung: Beispielpaper in einpflegen")
$ws.Range("D15").Value = "Eineleitung: Beispielpaper in einpflegen"

# Row 15: F15 replaced ("x" -> Webseite Banner text)
$ws.Range("F15").Value = "Webseite: Banner und Text für `"projekte`"…bild mobile lab o.ä."

# Row 16: F16 replaced ("x" -> Webseite Pilot invitation text)
$ws.Range("F16").Value = "Webseite: einladungstext für deinen Pilot schreiben, vielleicht Bild dazu, ähnlich dem Text für Mockdataparty"

# Update selection to reflect the last worked-on cell
$ws.Range("F30").Select() | Out-Null

$wb.Save()
